# Auto-generated edit script: refresh the cryptos price table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.295.94'
$ws.Range('E2').Value = '  +0.70%  '

$ws.Range('D3').Value = '1.664.28'
$ws.Range('E3').Value = '  +0.72%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.72'
$ws.Range('E5').Value = '  +0.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5329'
$ws.Range('E6').Value = '  +0.85%  '

$ws.Range('E8').Value = '  +1.27%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06363'
$ws.Range('E9').Value = '  +0.77%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.56'
$ws.Range('E10').Value = '  +0.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07854'
$ws.Range('E11').Value = '  +1.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.555'
$ws.Range('E12').Value = '  +1.65%  '

$ws.Range('D13').Value = '1.667.62'
$ws.Range('E13').Value = '  +0.23%  '

$ws.Range('D14').Value = '1.893.22'
$ws.Range('E14').Value = '  +0.72%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5533'
$ws.Range('E15').Value = '  +1.39%  '

$ws.Range('D16').Value = '0.0₅8190'
$ws.Range('E16').Value = '  +0.82%  '

$ws.Range('E17').Value = '  +0.82%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.010'
$ws.Range('E18').Value = '  +0.76%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.676'
$ws.Range('E19').Value = '  +3.05%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.41'
$ws.Range('E20').Value = '  -0.79%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.15'
$ws.Range('E21').Value = '  +1.23%  '

$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.061'
$ws.Range('E22').Value = '  +1.34%  '

$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.011'
$ws.Range('E23').Value = '  +0.69%  '

$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.09'
$ws.Range('E24').Value = '  +3.55%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1226'
$ws.Range('E25').Value = '  -1.12%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.244'
$ws.Range('E26').Value = '  -0.16%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.12'
$ws.Range('E27').Value = '  -0.12%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.476'
$ws.Range('E28').Value = '  +2.69%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05844'
$ws.Range('E29').Value = '  -1.53%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.279'
$ws.Range('E30').Value = '  +0.05%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.581'
$ws.Range('E31').Value = '  +2.22%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.306'
$ws.Range('E32').Value = '  +2.42%  '

$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.615'
$ws.Range('E33').Value = '  +4.57%  '

$ws.Range('B34').Value = 'MXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.825'
$ws.Range('E34').Value = '  +2.34%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9598'
$ws.Range('E35').Value = '  +1.78%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.434'
$ws.Range('E36').Value = '  +0.92%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5818'
$ws.Range('E37').Value = '  +3.52%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01611'
$ws.Range('E38').Value = '  +0.44%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.883'
$ws.Range('E39').Value = '  +0.59%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8548'
$ws.Range('E40').Value = '  +0.97%  '

$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.98'
$ws.Range('E41').Value = '  +4.00%  '

$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.010'
$ws.Range('E42').Value = '  +0.70%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.048.30'
$ws.Range('E43').Value = '  +3.97%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.805.08'
$ws.Range('E44').Value = '  +0.52%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.38'
$ws.Range('E45').Value = '  +1.04%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈105'
$ws.Range('E46').Value = '  -1.22%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  +0.53%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4374'
$ws.Range('E48').Value = '  +2.00%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.950'
$ws.Range('E49').Value = '  +2.52%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05166'
$ws.Range('E50').Value = '  +0.31%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.445'
$ws.Range('E51').Value = '  -1.70%  '
